# Update Wnt5a-Ryk LR-pair sheet with new TPM-derived values.
#
# The "Sending cluster" = FAPs rows (originally rows 2-7) are replaced
# with newly computed "ECs" sending-cluster rows, and the original FAPs
# rows are appended afterwards (now rows 8-13) with I/J/O/P/S/T
# specificity columns recomputed against the new ECs+FAPs totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Ryk"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.125615
$ws.Range("H2").Value = 0.25123
$ws.Range("I2").Value = 0.02647478672532295
$ws.Range("J2").Value = 0.01780700335556722
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.3291175
$ws.Range("N2").Value = 24.658235
$ws.Range("O2").Value = 0.1991607983368005
$ws.Range("P2").Value = 0.1614453197874725
$ws.Range("Q2").Value = 1.5487220947625
$ws.Range("R2").Value = 6.19488837905
$ws.Range("S2").Value = 0.005272739660011848
$ws.Range("T2").Value = 0.002874857351196145

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Ryk"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.125615
$ws.Range("H3").Value = 0.25123
$ws.Range("I3").Value = 0.02647478672532295
$ws.Range("J3").Value = 0.01780700335556722
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 28.70072166666667
$ws.Range("N3").Value = 86.10216500000001
$ws.Range("O3").Value = 0.4636226915653649
$ws.Range("P3").Value = 0.563738303362699
$ws.Range("Q3").Value = 3.605241152158334
$ws.Range("R3").Value = 21.63144691295
$ws.Range("S3").Value = 0.01227431188021322
$ws.Range("T3").Value = 0.01003848985964135

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Ryk"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.125615
$ws.Range("H4").Value = 0.25123
$ws.Range("I4").Value = 0.02647478672532295
$ws.Range("J4").Value = 0.01780700335556722
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.073231
$ws.Range("N4").Value = 0.219693
$ws.Range("O4").Value = 0.001182951206605196
$ws.Range("P4").Value = 0.001438400057427841
$ws.Range("Q4").Value = 0.009198912065
$ws.Range("R4").Value = 0.05519347239000001
$ws.Range("S4").Value = 0.00003131838090133602
$ws.Range("T4").Value = 0.00002561359464926565

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Ryk"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.125615
$ws.Range("H5").Value = 0.25123
$ws.Range("I5").Value = 0.02647478672532295
$ws.Range("J5").Value = 0.01780700335556722
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 20.652629
$ws.Range("N5").Value = 41.30525799999999
$ws.Range("O5").Value = 0.3336162608064818
$ws.Range("P5").Value = 0.2704386825218454
$ws.Range("Q5").Value = 2.594279991835
$ws.Range("R5").Value = 10.37711996734
$ws.Range("S5").Value = 0.008832419352951324
$ws.Range("T5").Value = 0.004815702527141679

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Ryk"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.125615
$ws.Range("H6").Value = 0.25123
$ws.Range("I6").Value = 0.02647478672532295
$ws.Range("J6").Value = 0.01780700335556722
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.05032066666666666
$ws.Range("N6").Value = 0.150962
$ws.Range("O6").Value = 0.0008128646795825703
$ws.Range("P6").Value = 0.0009883963051595711
$ws.Range("Q6").Value = 0.006321030543333332
$ws.Range("R6").Value = 0.03792618326
$ws.Range("S6").Value = 0.00002152041902849653
$ws.Range("T6").Value = 0.00001760037632260672

# Row 7
$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Ryk"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.125615
$ws.Range("H7").Value = 0.25123
$ws.Range("I7").Value = 0.02647478672532295
$ws.Range("J7").Value = 0.01780700335556722
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.099323
$ws.Range("N7").Value = 0.297969
$ws.Range("O7").Value = 0.001604433405165134
$ws.Range("P7").Value = 0.001950897965395876
$ws.Range("Q7").Value = 0.012476458645
$ws.Range("R7").Value = 0.07485875187
$ws.Range("S7").Value = 0.00004247703221673059
$ws.Range("T7").Value = 0.00003473964661617363

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Ryk"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.619088000000001
$ws.Range("H8").Value = 13.857264
$ws.Range("I8").Value = 0.9735252132746771
$ws.Range("J8").Value = 0.9821929966444328
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.3291175
$ws.Range("N8").Value = 24.658235
$ws.Range("O8").Value = 0.1991607983368005
$ws.Range("P8").Value = 0.1614453197874725
$ws.Range("Q8").Value = 56.94927869484
$ws.Range("R8").Value = 341.69567216904
$ws.Range("S8").Value = 0.1938880586767887
$ws.Range("T8").Value = 0.1585704624362763

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Ryk"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.619088000000001
$ws.Range("H9").Value = 13.857264
$ws.Range("I9").Value = 0.9735252132746771
$ws.Range("J9").Value = 0.9821929966444328
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 28.70072166666667
$ws.Range("N9").Value = 86.10216500000001
$ws.Range("O9").Value = 0.4636226915653649
$ws.Range("P9").Value = 0.563738303362699
$ws.Range("Q9").Value = 132.57115904184
$ws.Range("R9").Value = 1193.14043137656
$ws.Range("S9").Value = 0.4513483796851517
$ws.Range("T9").Value = 0.5536998135030577

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Wnt5a"
$ws.Range("C10").Value = "Ryk"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.619088000000001
$ws.Range("H10").Value = 13.857264
$ws.Range("I10").Value = 0.9735252132746771
$ws.Range("J10").Value = 0.9821929966444328
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.073231
$ws.Range("N10").Value = 0.219693
$ws.Range("O10").Value = 0.001182951206605196
$ws.Range("P10").Value = 0.001438400057427841
$ws.Range("Q10").Value = 0.338260433328
$ws.Range("R10").Value = 3.044343899952
$ws.Range("S10").Value = 0.00115163282570386
$ws.Range("T10").Value = 0.001412786462778575

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Wnt5a"
$ws.Range("C11").Value = "Ryk"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.619088000000001
$ws.Range("H11").Value = 13.857264
$ws.Range("I11").Value = 0.9735252132746771
$ws.Range("J11").Value = 0.9821929966444328
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 20.652629
$ws.Range("N11").Value = 41.30525799999999
$ws.Range("O11").Value = 0.3336162608064818
$ws.Range("P11").Value = 0.2704386825218454
$ws.Range("Q11").Value = 95.396310782352
$ws.Range("R11").Value = 572.377864694112
$ws.Range("S11").Value = 0.3247838414535305
$ws.Range("T11").Value = 0.2656229799947037

# Row 12
$ws.Range("A12").Value = "FAPs"
$ws.Range("B12").Value = "Wnt5a"
$ws.Range("C12").Value = "Ryk"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 4.619088000000001
$ws.Range("H12").Value = 13.857264
$ws.Range("I12").Value = 0.9735252132746771
$ws.Range("J12").Value = 0.9821929966444328
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.05032066666666666
$ws.Range("N12").Value = 0.150962
$ws.Range("O12").Value = 0.0008128646795825703
$ws.Range("P12").Value = 0.0009883963051595711
$ws.Range("Q12").Value = 0.232435587552
$ws.Range("R12").Value = 2.091920287968
$ws.Range("S12").Value = 0.0007913442605540738
$ws.Range("T12").Value = 0.0009707959288369645

# Row 13
$ws.Range("A13").Value = "FAPs"
$ws.Range("B13").Value = "Wnt5a"
$ws.Range("C13").Value = "Ryk"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 4.619088000000001
$ws.Range("H13").Value = 13.857264
$ws.Range("I13").Value = 0.9735252132746771
$ws.Range("J13").Value = 0.9821929966444328
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.099323
$ws.Range("N13").Value = 0.297969
$ws.Range("O13").Value = 0.001604433405165134
$ws.Range("P13").Value = 0.001950897965395876
$ws.Range("Q13").Value = 0.458781677424
$ws.Range("R13").Value = 4.129035096816
$ws.Range("S13").Value = 0.001561956372948403
$ws.Range("T13").Value = 0.001916158318779703
